$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.482.42"
$ws.Range("E2").Value = "  -3.12%  "
$ws.Range("D3").Value = "2.226.66"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "296.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.32%  "
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.15%  "
$ws.Range("E11").Value = "  -3.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.10%  "
$ws.Range("E14").Value = "  +9.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.104"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("D17").Value = "2.559.76"
$ws.Range("D18").Value = "2.245.59"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "42.386.06"
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.54%  "
$ws.Range("E21").Value = "  -4.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +22.19%  "
$ws.Range("E24").Value = "  -2.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "229.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.13%  "
$ws.Range("E31").Value = "  -3.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0903"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.95%  "
$ws.Range("E35").Value = "  +12.70%  "
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("E39").Value = "  -2.58%  "
$ws.Range("E40").Value = "  -3.26%  "
$ws.Range("E41").Value = "  -5.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.235"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.63%  "
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  -5.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  +6.76%  "
